$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 header ("Age (in years)" block) becomes "By age (in years)" style wording
$ws.Range("A18").Value = "Жаш курагы боюнча (жылдарда)"
$ws.Range("B18").Value = "По возрасту (в годах)"
$ws.Range("C18").Value = "By age (in years) "

# Row 28 header ("Education" block) becomes "By education" style wording
$ws.Range("A28").Value = "Билими боюнча"
$ws.Range("B28").Value = "По образованию"
$ws.Range("C28").Value = "By education"

# Clear the lingering D4 selection left over from editing, restoring A1 as the selection
$ws.Range("A1").Select()
